# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet
# for the currently-populated data rows (2-8), changing the timestamp from
# 2025-12-06 12:34:19 to 2025-12-06 12:44:18 (new append run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-06 12:34:19"
$newTimestamp = "2025-12-06 12:44:18"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
